$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.151.27"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "3.115.50"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "527.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.113.29"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.67%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D13").Value = "3.650.36"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.46%  "
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "58.199.07"
$ws.Range("E17").Value = "  -1.38%  "
$ws.Range("D18").Value = "3.123.61"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.59%  "
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "343.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.517"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("E26").Value = "  -0.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("D28").Value = "0.0₃0931"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.09%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("E34").Value = "  -1.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0670"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").Value = "3.156.34"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").Value = "2.272.75"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
